$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_9_1_23"
$ws.Range("B2").Value = 0.08636133185833772
$ws.Range("C2").Value = -0.3525509537251839
$ws.Range("D2").Value = -4.093905486246685
$ws.Range("E2").Value = -1.540784943869661
$ws.Range("F2").Value = 1.011128783226013
$ws.Range("G2").Value = 1.268031477928162
$ws.Range("H2").Value = 5.086646556854248
$ws.Range("I2").Value = 3.065027236938477

$ws.Range("A3").Value = "model_9_1_22"
$ws.Range("B3").Value = 0.1060899788836089
$ws.Range("C3").Value = -0.3138123153766257
$ws.Range("D3").Value = -3.968876390846575
$ws.Range("E3").Value = -1.476140228037995
$ws.Range("F3").Value = 0.9892950057983398
$ws.Range("G3").Value = 1.231713652610779
$ws.Range("H3").Value = 4.961795806884766
$ws.Range("I3").Value = 2.987044095993042

$ws.Range("A4").Value = "model_9_1_21"
$ws.Range("B4").Value = 0.1129634200882964
$ws.Range("C4").Value = -0.3421925102397121
$ws.Range("D4").Value = -3.872997361972358
$ws.Range("E4").Value = -1.450470072272106
$ws.Range("F4").Value = 0.9816880226135254
$ws.Range("G4").Value = 1.258320331573486
$ws.Range("H4").Value = 4.866053581237793
$ws.Range("I4").Value = 2.956077575683594

$ws.Range("A5").Value = "model_9_1_20"
$ws.Range("B5").Value = 0.1201266205805493
$ws.Range("C5").Value = -0.3303032044201244
$ws.Range("D5").Value = -3.828312193905566
$ws.Range("E5").Value = -1.428170942948769
$ws.Range("F5").Value = 0.9737604856491089
$ws.Range("G5").Value = 1.24717378616333
$ws.Range("H5").Value = 4.821431636810303
$ws.Range("I5").Value = 2.929177284240723

$ws.Range("A6").Value = "model_9_1_19"
$ws.Range("B6").Value = 0.1302341091797119
$ws.Range("C6").Value = -0.3020273159026903
$ws.Range("D6").Value = -3.777911964529959
$ws.Range("E6").Value = -1.396904295119169
$ws.Range("F6").Value = 0.9625745415687561
$ws.Range("G6").Value = 1.220664978027344
$ws.Range("H6").Value = 4.771103382110596
$ws.Range("I6").Value = 2.891459465026855

$ws.Range("A7").Value = "model_9_1_18"
$ws.Range("B7").Value = 0.1458310268250544
$ws.Range("C7").Value = -0.2621672434136835
$ws.Range("D7").Value = -3.697763136843732
$ws.Range("E7").Value = -1.349283211603367
$ws.Range("F7").Value = 0.9453133344650269
$ws.Range("G7").Value = 1.183295726776123
$ws.Range("H7").Value = 4.691068649291992
$ws.Range("I7").Value = 2.834012508392334

$ws.Range("A8").Value = "model_9_1_17"
$ws.Range("B8").Value = 0.1976626377650654
$ws.Range("C8").Value = 0.1020939470163857
$ws.Range("D8").Value = -3.677544458248595
$ws.Range("E8").Value = -1.191540544976411
$ws.Range("F8").Value = 0.8879510164260864
$ws.Range("G8").Value = 0.841796875
$ws.Range("H8").Value = 4.670878887176514
$ws.Range("I8").Value = 2.643723011016846

$ws.Range("A9").Value = "model_9_1_16"
$ws.Range("B9").Value = 0.3018288805323057
$ws.Range("C9").Value = 0.6135872599225978
$ws.Range("D9").Value = -3.405223518211733
$ws.Range("E9").Value = -0.875008436344245
$ws.Range("F9").Value = 0.7726697325706482
$ws.Range("G9").Value = 0.3622662127017975
$ws.Range("H9").Value = 4.398945808410645
$ws.Range("I9").Value = 2.26188063621521

$ws.Range("A10").Value = "model_9_1_15"
$ws.Range("B10").Value = 0.3469832922160264
$ws.Range("C10").Value = 0.6221419444280171
$ws.Range("D10").Value = -3.053916116114168
$ws.Range("E10").Value = -0.7346390568067662
$ws.Range("F10").Value = 0.7226970791816711
$ws.Range("G10").Value = 0.3542461097240448
$ws.Range("H10").Value = 4.048139095306396
$ws.Range("I10").Value = 2.092548370361328

$ws.Range("A11").Value = "model_9_1_14"
$ws.Range("B11").Value = 0.3654764275295174
$ws.Range("C11").Value = 0.5778903049688469
$ws.Range("D11").Value = -2.856898957032618
$ws.Range("E11").Value = -0.6760993977973424
$ws.Range("F11").Value = 0.7022305727005005
$ws.Range("G11").Value = 0.3957325220108032
$ws.Range("H11").Value = 3.851402759552002
$ws.Range("I11").Value = 2.02193021774292

$ws.Range("A12").Value = "model_9_1_9"
$ws.Range("B12").Value = 0.3702726609022364
$ws.Range("C12").Value = 0.5003217078588784
$ws.Range("D12").Value = -2.716286527299975
$ws.Range("E12").Value = -0.6532397590358159
$ws.Range("F12").Value = 0.6969226002693176
$ws.Range("G12").Value = 0.4684538841247559
$ws.Range("H12").Value = 3.71099066734314
$ws.Range("I12").Value = 1.994353771209717

$ws.Range("A13").Value = "model_9_1_13"
$ws.Range("B13").Value = 0.3708915448956226
$ws.Range("C13").Value = 0.5646305062969917
$ws.Range("D13").Value = -2.799782382332491
$ws.Range("E13").Value = -0.6593056015267202
$ws.Range("F13").Value = 0.6962376236915588
$ws.Range("G13").Value = 0.4081636965274811
$ws.Range("H13").Value = 3.794367551803589
$ws.Range("I13").Value = 2.001671552658081

$ws.Range("A14").Value = "model_9_1_12"
$ws.Range("B14").Value = 0.3741886414446604
$ws.Range("C14").Value = 0.5622628173183266
$ws.Range("D14").Value = -2.769544985536293
$ws.Range("E14").Value = -0.6485008989233352
$ws.Range("F14").Value = 0.692588746547699
$ws.Range("G14").Value = 0.4103834331035614
$ws.Range("H14").Value = 3.76417350769043
$ws.Range("I14").Value = 1.988637328147888

$ws.Range("A15").Value = "model_9_1_10"
$ws.Range("B15").Value = 0.3751661924355165
$ws.Range("C15").Value = 0.5586939113433826
$ws.Range("D15").Value = -2.752904630206094
$ws.Range("E15").Value = -0.6434868588391796
$ws.Range("F15").Value = 0.6915069222450256
$ws.Range("G15").Value = 0.413729339838028
$ws.Range("H15").Value = 3.747556924819946
$ws.Range("I15").Value = 1.982588648796082

$ws.Range("A16").Value = "model_9_1_11"
$ws.Range("B16").Value = 0.3758764504027238
$ws.Range("C16").Value = 0.5639893396346256
$ws.Range("D16").Value = -2.75572327735555
$ws.Range("E16").Value = -0.6424061881864607
$ws.Range("F16").Value = 0.6907208561897278
$ws.Range("G16").Value = 0.4087648093700409
$ws.Range("H16").Value = 3.750371217727661
$ws.Range("I16").Value = 1.981285095214844

$ws.Range("A17").Value = "model_9_1_8"
$ws.Range("B17").Value = 0.4177313757030275
$ws.Range("C17").Value = 0.6635831983334229
$ws.Range("D17").Value = -2.501775029253346
$ws.Range("E17").Value = -0.5025068348768555
$ws.Range("F17").Value = 0.6443997621536255
$ws.Range("G17").Value = 0.3153944909572601
$ws.Range("H17").Value = 3.496785163879395
$ws.Range("I17").Value = 1.812520265579224

$ws.Range("A18").Value = "model_9_1_7"
$ws.Range("B18").Value = 0.4291419577586616
$ws.Range("C18").Value = 0.7047443793887194
$ws.Range("D18").Value = -2.431668380783159
$ws.Range("E18").Value = -0.4582612808440769
$ws.Range("F18").Value = 0.6317716240882874
$ws.Range("G18").Value = 0.2768054306507111
$ws.Range("H18").Value = 3.426778316497803
$ws.Range("I18").Value = 1.759145379066467

$ws.Range("A19").Value = "model_9_1_6"
$ws.Range("B19").Value = 0.4350130275806344
$ws.Range("C19").Value = 0.7231283338741678
$ws.Range("D19").Value = -2.369046674237544
$ws.Range("E19").Value = -0.4263034141046473
$ws.Range("F19").Value = 0.625274121761322
$ws.Range("G19").Value = 0.2595702409744263
$ws.Range("H19").Value = 3.364245891571045
$ws.Range("I19").Value = 1.720593571662903

$ws.Range("A20").Value = "model_9_1_24"
$ws.Range("B20").Value = 0.4524559288649399
$ws.Range("C20").Value = -0.393361507491053
$ws.Range("D20").Value = -1.443112217113581
$ws.Range("E20").Value = -0.5249776675879214
$ws.Range("F20").Value = 0.6059699058532715
$ws.Range("G20").Value = 1.306291937828064
$ws.Range("H20").Value = 2.43963098526001
$ws.Range("I20").Value = 1.839627504348755

$ws.Range("A21").Value = "model_9_1_5"
$ws.Range("B21").Value = 0.4566363517473746
$ws.Range("C21").Value = 0.7479433297963354
$ws.Range("D21").Value = -2.187236408022287
$ws.Range("E21").Value = -0.3452711831059074
$ws.Range("F21").Value = 0.6013435125350952
$ws.Range("G21").Value = 0.2363059222698212
$ws.Range("H21").Value = 3.182694673538208
$ws.Range("I21").Value = 1.622841954231262

$ws.Range("A22").Value = "model_9_1_4"
$ws.Range("B22").Value = 0.5234496583640256
$ws.Range("C22").Value = 0.7537631714456523
$ws.Range("D22").Value = -1.607799699708335
$ws.Range("E22").Value = -0.1171612382495146
$ws.Range("F22").Value = 0.52740079164505
$ws.Range("G22").Value = 0.2308497428894043
$ws.Range("H22").Value = 2.60408353805542
$ws.Range("I22").Value = 1.347666025161743

$ws.Range("A23").Value = "model_9_1_3"
$ws.Range("B23").Value = 0.6272394327627293
$ws.Range("C23").Value = 0.7538633124691958
$ws.Range("D23").Value = -0.7305262825476537
$ws.Range("E23").Value = 0.2246155856860548
$ws.Range("F23").Value = 0.412536084651947
$ws.Range("G23").Value = 0.2307558655738831
$ws.Range("H23").Value = 1.728060245513916
$ws.Range("I23").Value = 0.9353700876235962

$ws.Range("A24").Value = "model_9_1_0"
$ws.Range("B24").Value = 0.6366181101798176
$ws.Range("C24").Value = 0.7676233475124017
$ws.Range("D24").Value = -0.4884726107533832
$ws.Range("E24").Value = 0.324567738887438
$ws.Range("F24").Value = 0.40215665102005
$ws.Range("G24").Value = 0.21785569190979
$ws.Range("H24").Value = 1.486351609230042
$ws.Range("I24").Value = 0.8147947788238525

$ws.Range("A25").Value = "model_9_1_1"
$ws.Range("B25").Value = 0.6385624802504575
$ws.Range("C25").Value = 0.7751767869358352
$ws.Range("D25").Value = -0.6023925339513891
$ws.Range("E25").Value = 0.2832986512887367
$ws.Range("F25").Value = 0.4000048041343689
$ws.Range("G25").Value = 0.2107742428779602
$ws.Range("H25").Value = 1.600109219551086
$ws.Range("I25").Value = 0.864578902721405

$ws.Range("A26").Value = "model_9_1_2"
$ws.Range("B26").Value = 0.6441613477637607
$ws.Range("C26").Value = 0.7708556754108795
$ws.Range("D26").Value = -0.590211714653214
$ws.Range("E26").Value = 0.2862656231255031
$ws.Range("F26").Value = 0.393808513879776
$ws.Range("G26").Value = 0.2148253321647644
$ws.Range("H26").Value = 1.587945699691772
$ws.Range("I26").Value = 0.8609997034072876
